# Updated via Streamlit Approval System
# Appends 5 new pending-approval rows (55-59) to the sheet, mirroring the
# existing row layout (columns A:AO).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 55 ----
$ws.Range("A55").Value = "WGG 02"
$ws.Range("B55").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("C55").Value = "15-01-2026"
$ws.Range("D55").Value = 286962
$ws.Range("E55").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("F55").Value = 34413429360
$ws.Range("G55").Value = "NEFT"
$ws.Range("H55").Value = "SBIN0003229"
$ws.Range("I55").Value = "AAAFW8862C"
$ws.Range("J55").Value = "32AAAFW8862C1Z9"
$ws.Range("L55").Value = "f6b19e08-2f51-45f7-b433-14be1e145835"
$ws.Range("U55").Value = "pending"
$ws.Range("V55").Value = 126000
$ws.Range("X55").Value = "Kolkata RPA_UNIQUE_ID : 48d30726-2a2a-4d62-bac1-04184b2abc77"
$ws.Range("Y55").Value = "Kolkata"
$ws.Range("Z55").Value = "PAYMENT"
$ws.Range("AA55").Value = "Payments@westernidc.com"
$ws.Range("AB55").Value = "ESTIMATION NOT MATCHED"
$ws.Range("AC55").Value = 0
$ws.Range("AD55").Value = 0
$ws.Range("AE55").Value = 0

# ---- Row 56 ----
$ws.Range("A56").Value = "WGG 02"
$ws.Range("B56").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("C56").Value = "15-01-2026"
$ws.Range("D56").Value = 286962
$ws.Range("E56").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("F56").Value = 34413429360
$ws.Range("G56").Value = "NEFT"
$ws.Range("H56").Value = "SBIN0003229"
$ws.Range("I56").Value = "AAAFW8862C"
$ws.Range("J56").Value = "32AAAFW8862C1Z9"
$ws.Range("L56").Value = "fd00fdd4-1fa3-4bc9-a26f-f2b9d22d1552"
$ws.Range("U56").Value = "pending"
$ws.Range("V56").Value = 66000
$ws.Range("X56").Value = "Income tax payamet  Hijas Sir 2024-25 RPA_UNIQUE_ID : ab1c1016-ed28-4a69-b600-e2c348ddce87"
$ws.Range("Y56").Value = "HO"
$ws.Range("Z56").Value = "PAYMENT"
$ws.Range("AA56").Value = "Payments@westernidc.com"
$ws.Range("AB56").Value = "ESTIMATION NOT MATCHED"
$ws.Range("AC56").Value = 0
$ws.Range("AD56").Value = 0
$ws.Range("AE56").Value = 0

# ---- Row 57 ----
$ws.Range("A57").Value = "WGG 02"
$ws.Range("B57").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("C57").Value = "15-01-2026"
$ws.Range("D57").Value = 286962
$ws.Range("E57").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("F57").Value = 34413429360
$ws.Range("G57").Value = "NEFT"
$ws.Range("H57").Value = "SBIN0003229"
$ws.Range("I57").Value = "AAAFW8862C"
$ws.Range("J57").Value = "32AAAFW8862C1Z9"
$ws.Range("L57").Value = "a9d5d8eb-98a8-47d5-b2eb-cf761ee76453"
$ws.Range("U57").Value = "pending"
$ws.Range("V57").Value = 66000
$ws.Range("X57").Value = "Income tax payamet  Hisham Sir 2024-25 RPA_UNIQUE_ID : e6530bda-5533-447d-b57c-414129b739d4"
$ws.Range("Y57").Value = "HO"
$ws.Range("Z57").Value = "PAYMENT"
$ws.Range("AA57").Value = "Payments@westernidc.com"
$ws.Range("AB57").Value = "ESTIMATION NOT MATCHED"
$ws.Range("AC57").Value = 0
$ws.Range("AD57").Value = 0
$ws.Range("AE57").Value = 0

# ---- Row 58 ----
$ws.Range("A58").Value = "WGG 02"
$ws.Range("B58").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("C58").Value = "15-01-2026"
$ws.Range("D58").Value = 286962
$ws.Range("E58").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("F58").Value = 34413429360
$ws.Range("G58").Value = "NEFT"
$ws.Range("H58").Value = "SBIN0003229"
$ws.Range("I58").Value = "AAAFW8862C"
$ws.Range("J58").Value = "32AAAFW8862C1Z9"
$ws.Range("L58").Value = "3e02b202-c072-460b-8784-389a4076b4dc"
$ws.Range("U58").Value = "pending"
$ws.Range("V58").Value = 30000
$ws.Range("X58").Value = "Income tax payamet Western 2024-25 RPA_UNIQUE_ID : 5ffcd7cd-1ef5-4889-ac17-fa65c60ff143"
$ws.Range("Y58").Value = "HO"
$ws.Range("Z58").Value = "PAYMENT"
$ws.Range("AA58").Value = "Payments@westernidc.com"
$ws.Range("AB58").Value = "ESTIMATION NOT MATCHED"
$ws.Range("AC58").Value = 0
$ws.Range("AD58").Value = 0
$ws.Range("AE58").Value = 0

# ---- Row 59 ----
$ws.Range("A59").Value = "WGG 02"
$ws.Range("B59").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("C59").Value = "15-01-2026"
$ws.Range("D59").Value = 286962
$ws.Range("E59").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("F59").Value = 34413429360
$ws.Range("G59").Value = "NEFT"
$ws.Range("H59").Value = "SBIN0003229"
$ws.Range("I59").Value = "AAAFW8862C"
$ws.Range("J59").Value = "32AAAFW8862C1Z9"
$ws.Range("L59").Value = "6665cd44-25e4-4005-a110-5cbb9f25d987"
$ws.Range("U59").Value = "pending"
$ws.Range("V59").Value = 38626
$ws.Range("X59").Value = "SIDBI Due RPA_UNIQUE_ID : 0f1272f9-af75-4c33-a61b-1f8ecf73bd0c"
$ws.Range("Y59").Value = "HO"
$ws.Range("Z59").Value = "PAYMENT"
$ws.Range("AA59").Value = "Payments@westernidc.com"
$ws.Range("AB59").Value = "ESTIMATION NOT MATCHED"
$ws.Range("AC59").Value = 0
$ws.Range("AD59").Value = 0
$ws.Range("AE59").Value = 0
